# Updated symbol list (Price / Volume(1h) columns) to match the refreshed
# coinranking.com snapshot. Cells D (Price) and E (Volume 1h) are stored as
# plain text (e.g. "300.97", "-0.03%"), so NumberFormat is forced to "@"
# (Text) before each write to keep Excel from re-interpreting the literal
# as a number/percentage.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "300.97"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "-0.03%"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "32.70"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "3.90%"
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "4.936"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "-2.94%"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07736"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "-1.54%"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "1.964"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "7.837"
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "0.40%"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.797"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "-0.60%"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.9207"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "-0.13%"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.1758"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "0.01%"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07794"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "2.86%"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.08593"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "-6.84%"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.03172"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "5.90%"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "-0.01%"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.001514"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "0.19%"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.005889"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "1.26%"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "3.460"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "-0.38%"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.153"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "-4.16%"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.1326"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "2.97%"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "4.96%"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.1992"
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "16.54%"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.04541"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "-1.28%"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001227"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "-1.76%"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004412"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "-1.34%"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "0.25%"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01708"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "-2.77%"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.04681"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "0.02%"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007618"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "7.18%"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1351"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "-0.75%"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.002344"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "7.15%"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.01144"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "17.00%"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.00006242"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "-0.28%"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "0.21%"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.8234"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "10.40%"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.00002102"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "0.21%"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0002002"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "0.21%"
